# StudentSem.xlsx edit
#
# Summary of changes (per commit "Adjust the weight and improve timetable exporter"):
#  1. Re-bucket the "Semester" values (column A) for the existing 400 data rows:
#       rows 127-201 : 1 -> 2
#       rows 252-376 : 2 -> 3
#       rows 377-401 : 2 -> 4
#  2. Append 100 new student rows (402-501), continuing the StudentID sequence
#     (step of 5) with Semester = 4 and Programme = "BCS".
#  3. Update the active selection to the newly appended block (A401:A501).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-bucket existing Semester values -------------------------------

$ws.Range("A127:A201").Value = 2
$ws.Range("A252:A376").Value = 3
$ws.Range("A377:A401").Value = 4

# --- 2. Append the 100 new rows (402-501) --------------------------------

$lastStudentId = 24005007
for ($row = 402; $row -le 501; $row++) {
    $lastStudentId = $lastStudentId + 5
    $ws.Cells.Item($row, 1).Value = 4
    $ws.Cells.Item($row, 2).Value = $lastStudentId
    $ws.Cells.Item($row, 3).Value = "BCS"
}

# --- 3. Update selection / view ------------------------------------------

$ws.Activate()
$ws.Range("A401:A501").Select()
